$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the quantity value in A6 (kategori_id) from 5 to 2
$ws.Range("A6").Value = 2

# Update the active cell selection to match the authored workbook state
$ws.Range("I13").Select()
